# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" sheets to match the refreshed scrape.
#
# Sheet "展览" (rows keyed by row number on that sheet):
#   F3: 38 -> 40
#   G4: 78 -> "已售罄" (sold out)
#   F5: 1798 -> 1801
#   F8: 157 -> 158
#   F9: 2206 -> 2223
#   F10: 4 -> 105
#   F11: 51 -> 56
#   F13: 1378 -> 1381
#   F14: 486 -> 487
#   F15: 27 -> 29
#   F23: 20 -> 43
#   F24: 21 -> 22
#   F25: 1194 -> 1197
#   F26: 9 -> 10
#   F27: 354 -> 359
#
# Sheet "全部类型" mirrors the same events but is offset by one extra row
# (it has an additional concert row inserted at row 6), so the same
# updates land on rows 3,4,5,9,10,11,12,14,15,16,24,25,26,27,28.

$wb = $excel.ActiveWorkbook

$updates1 = @{
    "F3"  = 40
    "F5"  = 1801
    "F8"  = 158
    "F9"  = 2223
    "F10" = 105
    "F11" = 56
    "F13" = 1381
    "F14" = 487
    "F15" = 29
    "F23" = 43
    "F24" = 22
    "F25" = 1197
    "F26" = 10
    "F27" = 359
}

$updates4 = @{
    "F3"  = 40
    "F5"  = 1801
    "F9"  = 158
    "F10" = 2223
    "F11" = 105
    "F12" = 56
    "F14" = 1381
    "F15" = 487
    "F16" = 29
    "F24" = 43
    "F25" = 22
    "F26" = 1197
    "F27" = 10
    "F28" = 359
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updates1.Keys) {
    $ws1.Range($addr).Value = $updates1[$addr]
}
$ws1.Range("G4").Value = "已售罄"

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updates4.Keys) {
    $ws4.Range($addr).Value = $updates4[$addr]
}
$ws4.Range("G4").Value = "已售罄"
